# Update NATMI LR-pair data (Wnt4-Fzd6) with refreshed TPM-based results.
# The new analysis adds the "Resolving-Mac" cluster to the sender/target cross-join,
# growing the table from a 4x3 (12-row) to a full 4x4 (16-row) cross-join, and refreshes
# every numeric column (expression / specificity / weight statistics) accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 16,20

# Row 2: sending ECs-ECs target
$data[0,0] = "ECs"
$data[0,1] = "Wnt4"
$data[0,2] = "Fzd6"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 1.497411
$data[0,7] = 4.492233
$data[0,8] = 0.3090503153498128
$data[0,9] = 0.3090503153498128
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 19.84402066666667
$data[0,13] = 59.532062
$data[0,14] = 0.8956779416773022
$data[0,15] = 0.8956779416773021
$data[0,16] = 29.71465483049399
$data[0,17] = 267.431893474446
$data[0,18] = 0.2768095503272414
$data[0,19] = 0.2768095503272414

# Row 3: sending ECs-FAPs target
$data[1,0] = "ECs"
$data[1,1] = "Wnt4"
$data[1,2] = "Fzd6"
$data[1,3] = "FAPs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 1.497411
$data[1,7] = 4.492233
$data[1,8] = 0.3090503153498128
$data[1,9] = 0.3090503153498128
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 1.745879666666666
$data[1,13] = 5.237639
$data[1,14] = 0.07880186845818919
$data[1,15] = 0.07880186845818919
$data[1,16] = 2.614299417543
$data[1,17] = 23.528694757887
$data[1,18] = 0.02435374229715783
$data[1,19] = 0.02435374229715783

# Row 4: sending ECs-MuSCs target
$data[2,0] = "ECs"
$data[2,1] = "Wnt4"
$data[2,2] = "Fzd6"
$data[2,3] = "MuSCs"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 1.497411
$data[2,7] = 4.492233
$data[2,8] = 0.3090503153498128
$data[2,9] = 0.3090503153498128
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 0.5556126666666666
$data[2,13] = 1.666838
$data[2,14] = 0.02507808362071368
$data[2,15] = 0.02507808362071367
$data[2,16] = 0.8319805188059999
$data[2,17] = 7.487824669253999
$data[2,18] = 0.007750389651350536
$data[2,19] = 0.007750389651350534

# Row 5: sending ECs-Resolving-Mac target
$data[3,0] = "ECs"
$data[3,1] = "Wnt4"
$data[3,2] = "Fzd6"
$data[3,3] = "Resolving-Mac"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 1.497411
$data[3,7] = 4.492233
$data[3,8] = 0.3090503153498128
$data[3,9] = 0.3090503153498128
$data[3,10] = 1
$data[3,11] = 0.3333333333333333
$data[3,12] = 0.009795
$data[3,13] = 0.029385
$data[3,14] = 0.0004421062437949407
$data[3,15] = 0.0004421062437949407
$data[3,16] = 0.014667140745
$data[3,17] = 0.132004266705
$data[3,18] = 0.0001366330740629476
$data[3,19] = 0.0001366330740629476

# Row 6: sending FAPs-ECs target
$data[4,0] = "FAPs"
$data[4,1] = "Wnt4"
$data[4,2] = "Fzd6"
$data[4,3] = "ECs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 1.589574
$data[4,7] = 4.768721999999999
$data[4,8] = 0.328071815935547
$data[4,9] = 0.3280718159355469
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 19.84402066666667
$data[4,13] = 59.532062
$data[4,14] = 0.8956779416773022
$data[4,15] = 0.8956779416773021
$data[4,16] = 31.543539307196
$data[4,17] = 283.8918537647639
$data[4,18] = 0.2938466888194854
$data[4,19] = 0.2938466888194854

# Row 7: sending FAPs-FAPs target
$data[5,0] = "FAPs"
$data[5,1] = "Wnt4"
$data[5,2] = "Fzd6"
$data[5,3] = "FAPs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 1.589574
$data[5,7] = 4.768721999999999
$data[5,8] = 0.328071815935547
$data[5,9] = 0.3280718159355469
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 1.745879666666666
$data[5,13] = 5.237639
$data[5,14] = 0.07880186845818919
$data[5,15] = 0.07880186845818919
$data[5,16] = 2.775204925261999
$data[5,17] = 24.97684432735799
$data[5,18] = 0.02585267208419223
$data[5,19] = 0.02585267208419222

# Row 8: sending FAPs-MuSCs target
$data[6,0] = "FAPs"
$data[6,1] = "Wnt4"
$data[6,2] = "Fzd6"
$data[6,3] = "MuSCs"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 1.589574
$data[6,7] = 4.768721999999999
$data[6,8] = 0.328071815935547
$data[6,9] = 0.3280718159355469
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 0.5556126666666666
$data[6,13] = 1.666838
$data[6,14] = 0.02507808362071368
$data[6,15] = 0.02507808362071367
$data[6,16] = 0.8831874490039998
$data[6,17] = 7.948687041035998
$data[6,18] = 0.008227412433631032
$data[6,19] = 0.00822741243363103

# Row 9: sending FAPs-Resolving-Mac target
$data[7,0] = "FAPs"
$data[7,1] = "Wnt4"
$data[7,2] = "Fzd6"
$data[7,3] = "Resolving-Mac"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 1.589574
$data[7,7] = 4.768721999999999
$data[7,8] = 0.328071815935547
$data[7,9] = 0.3280718159355469
$data[7,10] = 1
$data[7,11] = 0.3333333333333333
$data[7,12] = 0.009795
$data[7,13] = 0.029385
$data[7,14] = 0.0004421062437949407
$data[7,15] = 0.0004421062437949407
$data[7,16] = 0.01556987733
$data[7,17] = 0.14012889597
$data[7,18] = 0.0001450425982382499
$data[7,19] = 0.0001450425982382498

# Row 10: sending MuSCs-ECs target
$data[8,0] = "MuSCs"
$data[8,1] = "Wnt4"
$data[8,2] = "Fzd6"
$data[8,3] = "ECs"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 1.151
$data[8,7] = 3.453
$data[8,8] = 0.2375546279329019
$data[8,9] = 0.2375546279329019
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 19.84402066666667
$data[8,13] = 59.532062
$data[8,14] = 0.8956779416773022
$data[8,15] = 0.8956779416773021
$data[8,16] = 22.84046778733333
$data[8,17] = 205.564210086
$data[8,18] = 0.2127724401828589
$data[8,19] = 0.2127724401828589

# Row 11: sending MuSCs-FAPs target
$data[9,0] = "MuSCs"
$data[9,1] = "Wnt4"
$data[9,2] = "Fzd6"
$data[9,3] = "FAPs"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 1.151
$data[9,7] = 3.453
$data[9,8] = 0.2375546279329019
$data[9,9] = 0.2375546279329019
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 1.745879666666666
$data[9,13] = 5.237639
$data[9,14] = 0.07880186845818919
$data[9,15] = 0.07880186845818919
$data[9,16] = 2.009507496333333
$data[9,17] = 18.085567467
$data[9,18] = 0.01871974854200261
$data[9,19] = 0.01871974854200261

# Row 12: sending MuSCs-MuSCs target
$data[10,0] = "MuSCs"
$data[10,1] = "Wnt4"
$data[10,2] = "Fzd6"
$data[10,3] = "MuSCs"
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 1.151
$data[10,7] = 3.453
$data[10,8] = 0.2375546279329019
$data[10,9] = 0.2375546279329019
$data[10,10] = 3
$data[10,11] = 1
$data[10,12] = 0.5556126666666666
$data[10,13] = 1.666838
$data[10,14] = 0.02507808362071368
$data[10,15] = 0.02507808362071367
$data[10,16] = 0.6395101793333333
$data[10,17] = 5.755591613999999
$data[10,18] = 0.005957414823788839
$data[10,19] = 0.005957414823788837

# Row 13: sending MuSCs-Resolving-Mac target
$data[11,0] = "MuSCs"
$data[11,1] = "Wnt4"
$data[11,2] = "Fzd6"
$data[11,3] = "Resolving-Mac"
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 1.151
$data[11,7] = 3.453
$data[11,8] = 0.2375546279329019
$data[11,9] = 0.2375546279329019
$data[11,10] = 1
$data[11,11] = 0.3333333333333333
$data[11,12] = 0.009795
$data[11,13] = 0.029385
$data[11,14] = 0.0004421062437949407
$data[11,15] = 0.0004421062437949407
$data[11,16] = 0.011274045
$data[11,17] = 0.101466405
$data[11,18] = 0.00010502438425152
$data[11,19] = 0.0001050243842515199

# Row 14: sending Resolving-Mac-ECs target
$data[12,0] = "Resolving-Mac"
$data[12,1] = "Wnt4"
$data[12,2] = "Fzd6"
$data[12,3] = "ECs"
$data[12,4] = 3
$data[12,5] = 1
$data[12,6] = 0.6072163333333332
$data[12,7] = 1.821649
$data[12,8] = 0.1253232407817384
$data[12,9] = 0.1253232407817384
$data[12,10] = 3
$data[12,11] = 1
$data[12,12] = 19.84402066666667
$data[12,13] = 59.532062
$data[12,14] = 0.8956779416773022
$data[12,15] = 0.8956779416773021
$data[12,16] = 12.04961346780422
$data[12,17] = 108.446521210238
$data[12,18] = 0.1122492623477164
$data[12,19] = 0.1122492623477164

# Row 15: sending Resolving-Mac-FAPs target
$data[13,0] = "Resolving-Mac"
$data[13,1] = "Wnt4"
$data[13,2] = "Fzd6"
$data[13,3] = "FAPs"
$data[13,4] = 3
$data[13,5] = 1
$data[13,6] = 0.6072163333333332
$data[13,7] = 1.821649
$data[13,8] = 0.1253232407817384
$data[13,9] = 0.1253232407817384
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 1.745879666666666
$data[13,13] = 5.237639
$data[13,14] = 0.07880186845818919
$data[13,15] = 0.07880186845818919
$data[13,16] = 1.060126649634555
$data[13,17] = 9.541139846710999
$data[13,18] = 0.009875705534836521
$data[13,19] = 0.009875705534836521

# Row 16: sending Resolving-Mac-MuSCs target
$data[14,0] = "Resolving-Mac"
$data[14,1] = "Wnt4"
$data[14,2] = "Fzd6"
$data[14,3] = "MuSCs"
$data[14,4] = 3
$data[14,5] = 1
$data[14,6] = 0.6072163333333332
$data[14,7] = 1.821649
$data[14,8] = 0.1253232407817384
$data[14,9] = 0.1253232407817384
$data[14,10] = 3
$data[14,11] = 1
$data[14,12] = 0.5556126666666666
$data[14,13] = 1.666838
$data[14,14] = 0.02507808362071368
$data[14,15] = 0.02507808362071367
$data[14,16] = 0.3373770862068888
$data[14,17] = 3.036393775861999
$data[14,18] = 0.00314286671194327
$data[14,19] = 0.003142866711943269

# Row 17: sending Resolving-Mac-Resolving-Mac target
$data[15,0] = "Resolving-Mac"
$data[15,1] = "Wnt4"
$data[15,2] = "Fzd6"
$data[15,3] = "Resolving-Mac"
$data[15,4] = 3
$data[15,5] = 1
$data[15,6] = 0.6072163333333332
$data[15,7] = 1.821649
$data[15,8] = 0.1253232407817384
$data[15,9] = 0.1253232407817384
$data[15,10] = 1
$data[15,11] = 0.3333333333333333
$data[15,12] = 0.009795
$data[15,13] = 0.029385
$data[15,14] = 0.0004421062437949407
$data[15,15] = 0.0004421062437949407
$data[15,16] = 0.005947683984999999
$data[15,17] = 0.053529155865
$data[15,18] = 0.0000554061872422233
$data[15,19] = 0.0000554061872422233

# Write the full refreshed table back in one shot (sheet grows from A1:T13 to A1:T17).
$ws.Range("A2:T17").Value = $data
